$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5786272824987009
$ws.Range("C2").Value = 0.1424607078895122
$ws.Range("E2").Value = 0.5349065666312782
$ws.Range("F2").Value = 1.9386318224402
$ws.Range("G2").Value = 0.002397406590757517
$ws.Range("I2").Value = 0.3800386061133487
$ws.Range("J2").Value = 0.02881251915116678
$ws.Range("M2").Value = 0.5197229209234209
$ws.Range("N2").Value = 1.055473216310645
$ws.Range("O2").Value = 1.774413326933001
$ws.Range("B3").Value = 0.5075470041084031
$ws.Range("C3").Value = 0.1250525105573672
$ws.Range("E3").Value = 0.5307845632489077
$ws.Range("F3").Value = 1.924358166205451
$ws.Range("G3").Value = 0.002399921364743411
$ws.Range("I3").Value = 0.3845583701168991
$ws.Range("J3").Value = 0.02912978147959766
$ws.Range("M3").Value = 0.4893710185732161
$ws.Range("N3").Value = 1.063106144429121
$ws.Range("O3").Value = 1.773881814020569
$ws.Range("B4").Value = 0.463817675957813
$ws.Range("C4").Value = 0.114311521390249
$ws.Range("E4").Value = 0.5284867624169749
$ws.Range("F4").Value = 1.916738407699455
$ws.Range("G4").Value = 0.002401548064996573
$ws.Range("I4").Value = 0.3876088483800153
$ws.Range("J4").Value = 0.02933620189066133
$ws.Range("M4").Value = 0.470905199287877
$ws.Range("N4").Value = 1.068222459859513
$ws.Range("O4").Value = 1.774900766586669
$ws.Range("B5").Value = 0.4459771068897282
$ws.Range("C5").Value = 0.1099215316546633
$ws.Range("E5").Value = 0.5276090616776798
$ws.Range("F5").Value = 1.913920911138604
$ws.Range("G5").Value = 0.002402231795535941
$ws.Range("I5").Value = 0.3889210577113786
$ws.Range("J5").Value = 0.02942324358308568
$ws.Range("M5").Value = 0.4634234328729292
$ws.Range("N5").Value = 1.070415543999289
$ws.Range("O5").Value = 1.775653709435645
$ws.Range("B6").Value = 0.4430134860235171
$ws.Range("C6").Value = 0.1091918009407493
$ws.Range("E6").Value = 0.5274668652700782
$ws.Range("F6").Value = 1.913470435525412
$ws.Range("G6").Value = 0.002402346588863059
$ws.Range("I6").Value = 0.389143119970182
$ws.Range("J6").Value = 0.02943787342528559
$ws.Range("M6").Value = 0.4621837114612717
$ws.Range("N6").Value = 1.07078623863427
$ws.Range("O6").Value = 1.775799115108043
$ws.Range("B7").Value = 0.4635771534676678
$ws.Range("C7").Value = 0.1142523686218624
$ws.Range("E7").Value = 0.5284746878115172
$ws.Range("F7").Value = 1.91669924558397
$ws.Range("G7").Value = 0.002401557201585846
$ws.Range("I7").Value = 0.3876262656085103
$ws.Range("J7").Value = 0.02933736392533959
$ws.Range("M7").Value = 0.470804122152316
$ws.Range("N7").Value = 1.068251598552713
$ws.Range("O7").Value = 1.774909554424497
$ws.Range("B8").Value = 0.5541372074298749
$ws.Range("C8").Value = 0.1364693270177213
$ws.Range("E8").Value = 0.5334369421626732
$ws.Range("F8").Value = 1.933472671165319
$ws.Range("G8").Value = 0.002398256574872144
$ws.Range("I8").Value = 0.3815398061995197
$ws.Range("J8").Value = 0.02891950052751113
$ws.Range("M8").Value = 0.5092224321637175
$ws.Range("N8").Value = 1.058015931732101
$ws.Range("O8").Value = 1.773950490947499
$ws.Range("B9").Value = 0.7310071983799844
$ws.Range("C9").Value = 0.1796151617792532
$ws.Range("E9").Value = 0.5450165492370118
$ws.Range("F9").Value = 1.975455545941244
$ws.Range("G9").Value = 0.002392436756627445
$ws.Range("I9").Value = 0.3717941992304539
$ws.Range("J9").Value = 0.02819218542282531
$ws.Range("M9").Value = 0.5859006084301939
$ws.Range("N9").Value = 1.041348923741211
$ws.Range("O9").Value = 1.782773024147303
$ws.Range("B10").Value = 0.8604778505700779
$ws.Range("C10").Value = 0.211051401862818
$ws.Range("E10").Value = 0.5546512155846699
$ws.Range("F10").Value = 2.011863179507131
$ws.Range("G10").Value = 0.002388554858115927
$ws.Range("I10").Value = 0.3659765573218472
$ws.Range("J10").Value = 0.02771385173911511
$ws.Range("M10").Value = 0.6430431660340901
$ws.Range("N10").Value = 1.031174464289265
$ws.Range("O10").Value = 1.795824049989193
$ws.Range("B11").Value = 0.9192668243927073
$ws.Range("C11").Value = 0.2252943626271247
$ws.Range("E11").Value = 0.5592791694133865
$ws.Range("F11").Value = 2.029638851708867
$ws.Range("G11").Value = 0.002386873571152727
$ws.Range("I11").Value = 0.3636230663273317
$ws.Range("J11").Value = 0.02750838300262881
$ws.Range("M11").Value = 0.6692123571486093
$ws.Range("N11").Value = 1.026994540259658
$ws.Range("O11").Value = 1.8031974654576
$ws.Range("B12").Value = 0.9415122719584588
$ws.Range("C12").Value = 0.2306793665967746
$ws.Range("E12").Value = 0.5610668793326212
$ws.Range("F12").Value = 2.036544832242569
$ws.Range("G12").Value = 0.002386249014039853
$ws.Range("I12").Value = 0.3627741241417688
$ws.Range("J12").Value = 0.02743231963428805
$ws.Range("M12").Value = 0.6791467983182713
$ws.Range("N12").Value = 1.025476135749223
$ws.Range("O12").Value = 1.806196853004764
$ws.Range("B13").Value = 0.936722076060164
$ws.Range("C13").Value = 0.2295199910330155
$ws.Range("E13").Value = 0.5606802989318496
$ws.Range("F13").Value = 2.035049731946586
$ws.Range("G13").Value = 0.002386382985900438
$ws.Range("I13").Value = 0.3629550766161316
$ws.Range("J13").Value = 0.02744862372483414
$ws.Range("M13").Value = 0.6770061447553957
$ws.Range("N13").Value = 1.025800286002315
$ws.Range("O13").Value = 1.80554165362895
$ws.Range("B14").Value = 0.9210973112308807
$ws.Range("C14").Value = 0.2257375616875379
$ws.Range("E14").Value = 0.5594255400674299
$ws.Range("F14").Value = 2.030203508121161
$ws.Range("G14").Value = 0.0023868219461338
$ws.Range("I14").Value = 0.3635523752941907
$ws.Range("J14").Value = 0.02750209029186301
$ws.Range("M14").Value = 0.6700291760971027
$ws.Range("N14").Value = 1.026868328963836
$ws.Range("O14").Value = 1.80344006970401
$ws.Range("B15").Value = 0.9115244878323665
$ws.Range("C15").Value = 0.223419601481055
$ws.Range("E15").Value = 0.5586615474793177
$ws.Range("F15").Value = 2.027257815547387
$ws.Range("G15").Value = 0.002387092397127208
$ws.Range("I15").Value = 0.3639237475192623
$ws.Range("J15").Value = 0.02753506710015419
$ws.Range("M15").Value = 0.6657587910957972
$ws.Range("N15").Value = 1.027530926958221
$ws.Range("O15").Value = 1.802179797063133
$ws.Range("B16").Value = 0.8566335838942791
$ws.Range("C16").Value = 0.2101194140194593
$ws.Range("E16").Value = 0.5543536964204847
$ws.Range("F16").Value = 2.010725936227402
$ws.Range("G16").Value = 0.002388666431999198
$ws.Range("I16").Value = 0.3661362713068357
$ws.Range("J16").Value = 0.02772752357806763
$ws.Range("M16").Value = 0.6413364307494618
$ws.Range("N16").Value = 1.031456652504424
$ws.Range("O16").Value = 1.795371142675208
$ws.Range("B17").Value = 0.8229313446326501
$ws.Range("C17").Value = 0.2019452839121243
$ws.Range("E17").Value = 0.5517737154452647
$ws.Range("F17").Value = 2.000895150399629
$ws.Range("G17").Value = 0.002389653680916811
$ws.Range("I17").Value = 0.3675687249617887
$ws.Range("J17").Value = 0.0278486950263126
$ws.Range("M17").Value = 0.6263985827695961
$ws.Range("N17").Value = 1.033979789562601
$ws.Range("O17").Value = 1.791562642357604
$ws.Range("B18").Value = 0.8035366227433087
$ws.Range("C18").Value = 0.1972383375076276
$ws.Range("E18").Value = 0.5503128505722401
$ws.Range("F18").Value = 1.995354965759802
$ws.Range("G18").Value = 0.002390229487313256
$ws.Range("I18").Value = 0.3684202028614436
$ws.Range("J18").Value = 0.02791953119374124
$ws.Range("M18").Value = 0.6178231999372485
$ws.Range("N18").Value = 1.035473247831796
$ws.Range("O18").Value = 1.789507257702581
$ws.Range("B19").Value = 0.7969682060544869
$ws.Range("C19").Value = 0.1956437251301963
$ws.Range("E19").Value = 0.5498221906916285
$ws.Range("F19").Value = 1.993498766417645
$ws.Range("G19").Value = 0.002390425815575395
$ws.Range("I19").Value = 0.3687132291930446
$ws.Range("J19").Value = 0.02794371119029737
$ws.Range("M19").Value = 0.6149225637667683
$ws.Range("N19").Value = 1.035986159481283
$ws.Range("O19").Value = 1.788834533670439
$ws.Range("B20").Value = 0.826520057651237
$ws.Range("C20").Value = 0.2028159948102086
$ws.Range("E20").Value = 0.5520459712573427
$ws.Range("F20").Value = 2.001929831217197
$ws.Range("G20").Value = 0.002389547762201957
$ws.Range("I20").Value = 0.3674133835563254
$ws.Range("J20").Value = 0.02783567797744357
$ws.Range("M20").Value = 0.6279870402234167
$ws.Range("N20").Value = 1.033706828434404
$ws.Range("O20").Value = 1.79195406983763
$ws.Range("B21").Value = 0.9256871447238382
$ws.Range("C21").Value = 0.2268487852198575
$ws.Range("E21").Value = 0.5597931380465084
$ws.Range("F21").Value = 2.031622218132171
$ws.Range("G21").Value = 0.002386692684717488
$ws.Range("I21").Value = 0.3633757857557995
$ws.Range("J21").Value = 0.0274863385611841
$ws.Range("M21").Value = 0.6720778127063767
$ws.Range("N21").Value = 1.026552870295767
$ws.Range("O21").Value = 1.804051726156416
$ws.Range("B22").Value = 0.9904009251189336
$ws.Range("C22").Value = 0.2425059785018391
$ws.Range("E22").Value = 0.5650615383836595
$ws.Range("F22").Value = 2.052046398544249
$ws.Range("G22").Value = 0.002384897285689586
$ws.Range("I22").Value = 0.3609834524092079
$ws.Range("J22").Value = 0.02726818533224051
$ws.Range("M22").Value = 0.701037682238109
$ws.Range("N22").Value = 1.022252942355934
$ws.Range("O22").Value = 1.81316648181604
$ws.Range("B23").Value = 0.955871275776019
$ws.Range("C23").Value = 0.2341540534965247
$ws.Range("E23").Value = 0.5622309332340549
$ws.Range("F23").Value = 2.041052372136278
$ws.Range("G23").Value = 0.002385849086366396
$ws.Range("I23").Value = 0.3622376864561971
$ws.Range("J23").Value = 0.02738368839325389
$ws.Range("M23").Value = 0.6855682084812855
$ws.Range("N23").Value = 1.024513543034352
$ws.Range("O23").Value = 1.808190988086182
$ws.Range("B24").Value = 0.8248976591422661
$ws.Range("C24").Value = 0.2024223698532523
$ws.Range("E24").Value = 0.5519228146253283
$ws.Range("F24").Value = 2.001461704202029
$ws.Range("G24").Value = 0.00238959562245292
$ws.Range("I24").Value = 0.367483526380429
$ws.Range("J24").Value = 0.02784155933116761
$ws.Range("M24").Value = 0.6272688593785034
$ws.Range("N24").Value = 1.033830100606799
$ws.Range("O24").Value = 1.791776687580182
$ws.Range("B25").Value = 0.6832399881578795
$ws.Range("C25").Value = 0.1679888245533334
$ws.Range("E25").Value = 0.5416860357273805
$ws.Range("F25").Value = 1.963122708278235
$ws.Range("G25").Value = 0.002393941708945331
$ws.Range("I25").Value = 0.3741954307916160
$ws.Range("J25").Value = 0.02837909562382013
$ws.Range("M25").Value = 0.5650146899510773
$ws.Range("N25").Value = 1.045493791655517
$ws.Range("O25").Value = 1.779235873540273
